# Daily attendance processing - 2026-01-04 21:56:37
# Rotate the "Recorded By" (column G) comma-separated list of recorders so
# that the last recorder in the list is moved to the front (most-recent-first).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $txt = $cell.Text

    if ($txt -and $txt.Contains(",")) {
        $parts = $txt -split ",\s*"
        $count = $parts.Count

        if ($count -gt 1) {
            $rotated = @($parts[$count - 1]) + $parts[0..($count - 2)]
            $newVal = [string]::Join(", ", $rotated)
            $cell.Value2 = $newVal
        }
    }
}
